$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: rename "Image" column and add circle2-4 headers
$ws.Cells.Item(1, 3).Value = "circle1"
$ws.Cells.Item(1, 4).Value = "circle2"
$ws.Cells.Item(1, 5).Value = "circle3"
$ws.Cells.Item(1, 6).Value = "circle4"

# Replace numeric group values in column C with boolean circle-hit flags,
# and populate new boolean columns D:F for circle2/circle3/circle4
$ws.Cells.Item(2, 3).Value = $false
$ws.Cells.Item(2, 4).Value = $false
$ws.Cells.Item(2, 5).Value = $false
$ws.Cells.Item(2, 6).Value = $false
$ws.Cells.Item(3, 3).Value = $false
$ws.Cells.Item(3, 4).Value = $true
$ws.Cells.Item(3, 5).Value = $false
$ws.Cells.Item(3, 6).Value = $false
$ws.Cells.Item(4, 3).Value = $true
$ws.Cells.Item(4, 4).Value = $true
$ws.Cells.Item(4, 5).Value = $true
$ws.Cells.Item(4, 6).Value = $true
$ws.Cells.Item(5, 3).Value = $false
$ws.Cells.Item(5, 4).Value = $false
$ws.Cells.Item(5, 5).Value = $true
$ws.Cells.Item(5, 6).Value = $false
$ws.Cells.Item(6, 3).Value = $true
$ws.Cells.Item(6, 4).Value = $false
$ws.Cells.Item(6, 5).Value = $true
$ws.Cells.Item(6, 6).Value = $false
$ws.Cells.Item(7, 3).Value = $false
$ws.Cells.Item(7, 4).Value = $false
$ws.Cells.Item(7, 5).Value = $false
$ws.Cells.Item(7, 6).Value = $true
$ws.Cells.Item(8, 3).Value = $true
$ws.Cells.Item(8, 4).Value = $false
$ws.Cells.Item(8, 5).Value = $true
$ws.Cells.Item(8, 6).Value = $true
$ws.Cells.Item(9, 3).Value = $false
$ws.Cells.Item(9, 4).Value = $true
$ws.Cells.Item(9, 5).Value = $false
$ws.Cells.Item(9, 6).Value = $false
$ws.Cells.Item(10, 3).Value = $true
$ws.Cells.Item(10, 4).Value = $false
$ws.Cells.Item(10, 5).Value = $true
$ws.Cells.Item(10, 6).Value = $true
$ws.Cells.Item(11, 3).Value = $false
$ws.Cells.Item(11, 4).Value = $true
$ws.Cells.Item(11, 5).Value = $true
$ws.Cells.Item(11, 6).Value = $false
$ws.Cells.Item(12, 3).Value = $false
$ws.Cells.Item(12, 4).Value = $false
$ws.Cells.Item(12, 5).Value = $false
$ws.Cells.Item(12, 6).Value = $false
$ws.Cells.Item(13, 3).Value = $false
$ws.Cells.Item(13, 4).Value = $true
$ws.Cells.Item(13, 5).Value = $false
$ws.Cells.Item(13, 6).Value = $false
$ws.Cells.Item(14, 3).Value = $true
$ws.Cells.Item(14, 4).Value = $false
$ws.Cells.Item(14, 5).Value = $true
$ws.Cells.Item(14, 6).Value = $true
$ws.Cells.Item(15, 3).Value = $false
$ws.Cells.Item(15, 4).Value = $true
$ws.Cells.Item(15, 5).Value = $false
$ws.Cells.Item(15, 6).Value = $false
$ws.Cells.Item(16, 3).Value = $true
$ws.Cells.Item(16, 4).Value = $false
$ws.Cells.Item(16, 5).Value = $true
$ws.Cells.Item(16, 6).Value = $true
$ws.Cells.Item(17, 3).Value = $false
$ws.Cells.Item(17, 4).Value = $true
$ws.Cells.Item(17, 5).Value = $false
$ws.Cells.Item(17, 6).Value = $false
$ws.Cells.Item(18, 3).Value = $true
$ws.Cells.Item(18, 4).Value = $false
$ws.Cells.Item(18, 5).Value = $true
$ws.Cells.Item(18, 6).Value = $true
$ws.Cells.Item(19, 3).Value = $false
$ws.Cells.Item(19, 4).Value = $true
$ws.Cells.Item(19, 5).Value = $false
$ws.Cells.Item(19, 6).Value = $false
$ws.Cells.Item(20, 3).Value = $true
$ws.Cells.Item(20, 4).Value = $false
$ws.Cells.Item(20, 5).Value = $true
$ws.Cells.Item(20, 6).Value = $true
$ws.Cells.Item(21, 3).Value = $false
$ws.Cells.Item(21, 4).Value = $true
$ws.Cells.Item(21, 5).Value = $false
$ws.Cells.Item(21, 6).Value = $false
$ws.Cells.Item(22, 3).Value = $true
$ws.Cells.Item(22, 4).Value = $false
$ws.Cells.Item(22, 5).Value = $true
$ws.Cells.Item(22, 6).Value = $true
$ws.Cells.Item(23, 3).Value = $false
$ws.Cells.Item(23, 4).Value = $true
$ws.Cells.Item(23, 5).Value = $false
$ws.Cells.Item(23, 6).Value = $false
$ws.Cells.Item(24, 3).Value = $true
$ws.Cells.Item(24, 4).Value = $false
$ws.Cells.Item(24, 5).Value = $true
$ws.Cells.Item(24, 6).Value = $true

# Page setup now carries an explicit (default) print orientation
$ws.PageSetup.Orientation = 1

# Move active selection to F7 to match the edited workbook state
$ws.Range("F7").Select()
